# Mesquite_Starter_INSTRUCTIONS_Windows.docx - Release of 3.70, build 940
#
# Two edits:
# 1) In the "PREREQUISITE: JAVA" paragraph, split the sentence so that
#    "...or higher, but we have tested..." becomes
#    "...or higher. We have tested..." and the run is split into three
#    runs: "...or higher" | ". W" | "e have tested...".
# 2) Merge the "This forces us..." paragraph, the following empty
#    paragraph, and the "The problems are usually..." paragraph into a
#    single paragraph separated by a plain space.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1
# ---------------------------------------------------------------------

# Replace ", but w" with ". W" (this collapses the sentence into one run
# first; we then re-split the run below so the saved XML has three runs,
# matching the target edit exactly).
$find1 = $d.Content
$find1.Find.Execute(", but w", $true, $false, $false, $false, $false, $true, 1, $false, ". W", 2) | Out-Null

# Work out the two split points: right after "or higher" and right after
# the newly inserted ". W".
$afterHigher = $d.Content
$afterHigher.Find.Execute("or higher", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos1 = $afterHigher.End

$afterDotW = $d.Content
$afterDotW.Find.Execute(". W", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos2 = $afterDotW.End

# Forcing a run split: briefly placing (and removing) a bookmark at a
# collapsed range boundary causes the run containing that point to be
# broken into two independent runs, even after the bookmark is removed.
$d.Bookmarks.Add("zzSplit1", $d.Range($splitPos1, $splitPos1)) | Out-Null
$d.Bookmarks.Add("zzSplit2", $d.Range($splitPos2, $splitPos2)) | Out-Null
$d.Bookmarks.Item("zzSplit1").Delete()
$d.Bookmarks.Item("zzSplit2").Delete()

# ---------------------------------------------------------------------
# Edit 2
# ---------------------------------------------------------------------

# Locate the three paragraphs involved by their text.
$paraForces = $null
$paraEmpty = $null
$paraProblems = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*This forces us to maintain multiple versions*") {
        $paraForces = $i
    } elseif (($paraForces -ne $null) -and ($paraEmpty -eq $null) -and ($i -eq ($paraForces + 1))) {
        $paraEmpty = $i
    } elseif (($paraEmpty -ne $null) -and ($paraProblems -eq $null) -and ($i -eq ($paraEmpty + 1))) {
        $paraProblems = $i
    }
}

# Insert the separating space FIRST, while "This forces us..." is still
# its own paragraph (right before its own paragraph mark). Doing this
# before any paragraph marks are deleted creates a clean, dedicated run
# for the space, and avoids corrupting the xml:space handling of the
# runs being joined.
$pForces = $d.Paragraphs.Item($paraForces)
$endOfForces = $pForces.Range.End - 1
$d.Range($endOfForces, $endOfForces).InsertAfter(" ")

# Delete the paragraph mark ending the (still) empty paragraph. This
# merges it (no content) onto the start of the "The problems..."
# paragraph.
$pEmpty = $d.Paragraphs.Item($paraEmpty)
$markEmpty = $d.Range($pEmpty.Range.End - 1, $pEmpty.Range.End)
$markEmpty.Delete()

# Delete the paragraph mark ending the "This forces us..." paragraph
# (now ending with the space we inserted above). This merges it with
# the paragraph that used to be empty (now carrying "The problems..."),
# yielding one single paragraph: "...Java 8). The problems...start."
$pForces2 = $d.Paragraphs.Item($paraForces)
$markForces = $d.Range($pForces2.Range.End - 1, $pForces2.Range.End)
$markForces.Delete()
